# Apply crypto price/volume updates per commit "Updated cryptos list on Sat Apr 22 14:36:00 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 40/41 swap: FraxShare <-> Frax (name, link, price, change all swap) ---
$ws.Range("B40").Value = "Frax"
$ws.Range("C40").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"

# --- Price column (D): force text format so values like "27.367.17" / "47.10" / "8.020"
#     are stored verbatim as text instead of being auto-coerced into numbers ---
$priceCells = @("D2", "D3", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D35", "D36", "D38", "D40", "D41", "D42", "D43", "D45", "D46", "D48", "D49", "D51")
foreach ($cell in $priceCells) { $ws.Range($cell).NumberFormat = "@" }

$ws.Range("D2").Value = "27.367.17"
$ws.Range("D3").Value = "1.859.02"
$ws.Range("D5").Value = "330.23"
$ws.Range("D7").Value = "0.4699"
$ws.Range("D8").Value = "0.3958"
$ws.Range("D9").Value = "47.10"
$ws.Range("D10").Value = "0.08003"
$ws.Range("D11").Value = "1.015"
$ws.Range("D12").Value = "21.52"
$ws.Range("D13").Value = "1.854.23"
$ws.Range("D14").Value = "5.948"
$ws.Range("D15").Value = "7.145"
$ws.Range("D16").Value = "1.004"
$ws.Range("D18").Value = "0.00001036"
$ws.Range("D19").Value = "0.06555"
$ws.Range("D20").Value = "17.28"
$ws.Range("D21").Value = "1.001"
$ws.Range("D22").Value = "5.491"
$ws.Range("D23").Value = "27.355.36"
$ws.Range("D24").Value = "10.93"
$ws.Range("D25").Value = "2.298"
$ws.Range("D26").Value = "2.060.99"
$ws.Range("D27").Value = "20.43"
$ws.Range("D28").Value = "153.80"
$ws.Range("D29").Value = "2.071"
$ws.Range("D30").Value = "5.487"
$ws.Range("D31").Value = "122.18"
$ws.Range("D32").Value = "0.09509"
$ws.Range("D35").Value = "3.589"
$ws.Range("D36").Value = "5.279"
$ws.Range("D38").Value = "0.02230"
$ws.Range("D40").Value = "1.001"
$ws.Range("D41").Value = "8.020"
$ws.Range("D42").Value = "0.5939"
$ws.Range("D43").Value = "0.1899"
$ws.Range("D45").Value = "1.280"
$ws.Range("D46").Value = "0.5656"
$ws.Range("D48").Value = "3.429"
$ws.Range("D49").Value = "1.927"
$ws.Range("D51").Value = "109.44"

foreach ($cell in $priceCells) { $ws.Range($cell).Style = "Normal" }

# --- Volume(1h) column (E): plain text assignment (percent strings never parse as numbers) ---
$ws.Range("E2").Value = "  -2.87%  "
$ws.Range("E3").Value = "  -3.04%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("E5").Value = "  +0.77%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("E8").Value = "  -1.52%  "
$ws.Range("E9").Value = "  -10.99%  "
$ws.Range("E10").Value = "  -4.89%  "
$ws.Range("E11").Value = "  -3.00%  "
$ws.Range("E12").Value = "  -3.01%  "
$ws.Range("E13").Value = "  -2.29%  "
$ws.Range("E14").Value = "  -2.16%  "
$ws.Range("E15").Value = "  -4.09%  "
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("E17").Value = "  -3.76%  "
$ws.Range("E18").Value = "  -2.67%  "
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("E20").Value = "  -3.77%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("E22").Value = "  -4.48%  "
$ws.Range("E23").Value = "  -2.90%  "
$ws.Range("E24").Value = "  -2.43%  "
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("E26").Value = "  -3.23%  "
$ws.Range("E27").Value = "  +1.87%  "
$ws.Range("E28").Value = "  +0.28%  "
$ws.Range("E29").Value = "  -3.05%  "
$ws.Range("E30").Value = "  -5.05%  "
$ws.Range("E31").Value = "  -1.13%  "
$ws.Range("E32").Value = "  -1.69%  "
$ws.Range("E33").Value = "  -2.89%  "
$ws.Range("E34").Value = "  +0.38%  "
$ws.Range("E35").Value = "  -1.53%  "
$ws.Range("E36").Value = "  -4.88%  "
$ws.Range("E37").Value = "  -2.08%  "
$ws.Range("E38").Value = "  -3.08%  "
$ws.Range("E39").Value = "  -3.85%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("E41").Value = "  -9.44%  "
$ws.Range("E42").Value = "  -3.83%  "
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("E45").Value = "  -2.26%  "
$ws.Range("E46").Value = "  -3.50%  "
$ws.Range("E47").Value = "  -5.52%  "
$ws.Range("E48").Value = "  -0.24%  "
$ws.Range("E49").Value = "  -4.95%  "
$ws.Range("E51").Value = "  -1.73%  "
